# "Added New Mac-Address and Document Types"
# Appends one new test-data row (row 33) to the master-reg_center_user sheet,
# mirroring the existing rows' pattern (regcntr_id, usr_id, lang_code,
# is_active, cr_by, cr_dtimes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_user")

$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"

# Mirror the author's final on-screen state: scrolled down with C31 selected.
$ws.Range("C31").Select()
